$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume updates (GitHub Actions data refresh)
# Force text format so numeric-looking strings (e.g. "547.34") are not
# auto-converted to numbers, matching the original inlineStr text cells.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.329.54'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +4.25%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.347.89'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +2.79%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '547.34'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '132.57'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +1.21%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +1.41%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.345.81'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.79%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +2.49%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.79%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.17%  '
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '23.93'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +2.27%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.762.18'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +2.82%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '60.266.84'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +4.18%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000134'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.08%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.351.49'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +2.98%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.68'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +7.33%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '314.67'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.55'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.96%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.172'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.73%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.89'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.76%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.36'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +8.21%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +3.15%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '171.65'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +11.91%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0₃0730'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +2.21%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +4.28%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +15.32%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.383'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +1.47%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '18.04'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +2.00%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.18'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +8.25%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '326.23'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +14.88%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +4.00%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '38.12'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.92%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '141.59'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.71%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +1.84%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0951'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.57%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '19.51'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +9.03%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +1.24%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.562'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +2.18%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0215'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +3.38%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0₆0212'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +16.32%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '11.03'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.13%  '
